$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Assigned but not busy"
